# Apply the bill-of-quantities update described by the commit diff.
# Rows 8-12 get new item data; rows 14 & 16 get the recalculated grand totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 : "P. point" / Medium point (up to 6 mtr.) ---
$ws.Range("A8").Value = "P. point"
$ws.Range("C8").Value = 55
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3"
$ws.Range("E8").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F8").Value = 472
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "25960.00"

# --- Row 9 : Long point (up to 10 mtr.) ---
$ws.Range("C9").Value = 24
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "4"
$ws.Range("E9").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F9").Value = 662
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "15888.00"

# --- Row 10 : "P. point" / On board ---
$ws.Range("A10").Value = "P. point"
$ws.Range("C10").Value = 26
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6"
$ws.Range("E10").Value = "On board"
$ws.Range("F10").Value = 136
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "3536.00"

# --- Row 11 : "%" / Add Tender Premium ---
$ws.Range("A11").Value = "%"
$ws.Range("C11").Value = 56
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "9"
$ws.Range("E11").Value = "Add Tender Premium "

# --- Row 12 : Grand Total quantity revised ---
$ws.Range("C12").Value = 17

# --- Row 14 : Grand Total Rs. ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "45384.00"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "45384.00"

# --- Row 16 : NET PAYABLE AMOUNT Rs. ---
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "45384.00"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "45384.00"
